$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in Saturday/Sunday hours for the week of 43178 (row 10)
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 6

# Move the active selection to I20, mirroring where the author ended up
$ws.Range("I20").Select()
